$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Balance")

# Rename the sheet; Excel automatically repoints every defined name that
# held a real cell reference (TotalActivo, TotalPasivoPatrimonio, etc.) to
# the new sheet name.
$ws.Name = "ESF"

# Restore the print area (re-assigning keeps it anchored to the renamed sheet).
$ws.PageSetup.PrintArea = '$B$2:$J$54'

# A handful of defined names only ever held a broken #REF! reference
# qualified by the old sheet name (Balance!#REF!). Re-point those at the
# new sheet name explicitly, the same way Excel keeps them anchored to the
# renamed sheet.
$refNames = "Año", "AÑOINICIAL", "CentroCos.", "Mes", "MESINICIAL", "No.Empresa", "NombMes", "Unid.MOneda", "UNIDAD"
foreach ($nm in $refNames) {
    $wb.Names.Item($nm).RefersTo = "=ESF!#REF!"
}
